# Updated cryptos list — refresh prices / 1h volume % and fix the
# Maker/Hedera row ordering (rows 47-48 swapped content).
#
# NOTE: several "Price" column values are plain digit-group strings like
# "523.75" or "0.569" that Excel would otherwise auto-coerce into numbers
# (losing the original text formatting / exact string form, since these
# are inline strings in the source workbook, not numeric cells). We force
# those assignments to stay text the same way a user typing into Excel
# would: a leading apostrophe. Values containing more than one "." (e.g.
# "59.065.46") are never parsed as numbers by Excel, so they don't need
# the apostrophe, but using it everywhere in the Price column is harmless
# and keeps this script uniform/robust.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $value
    # The leading apostrophe forces Excel to keep a numeric-looking string
    # as text, but it also stamps the cell with a "quote prefix" style.
    # Reset back to the default style so formatting matches the original
    # (unstyled) cell exactly - only the text content should change.
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-Text "D2" "59.065.46"
$ws.Range("E2").Value = "  -2.10%  "

# Row 3 - Ethereum
Set-Text "D3" "2.664.45"
$ws.Range("E3").Value = "  -0.55%  "

# Row 4 - TetherUSD (price unchanged)
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
Set-Text "D5" "523.75"
$ws.Range("E5").Value = "  +0.14%  "

# Row 6 - Solana
Set-Text "D6" "144.36"
$ws.Range("E6").Value = "  -1.09%  "

# Row 7 - USDC (price unchanged)
$ws.Range("E7").Value = "  +0.20%  "

# Row 8 - XRP
Set-Text "D8" "0.569"
$ws.Range("E8").Value = "  -1.07%  "

# Row 9 - Toncoin
Set-Text "D9" "6.99"
$ws.Range("E9").Value = "  +8.51%  "

# Row 10 - Dogecoin (price unchanged)
$ws.Range("E10").Value = "  -2.39%  "

# Row 11 - Cardano (price unchanged)
$ws.Range("E11").Value = "  -1.82%  "

# Row 12 - TRON (price unchanged)
$ws.Range("E12").Value = "  +1.38%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-Text "D13" "3.132.87"
$ws.Range("E13").Value = "  -0.62%  "

# Row 14 - WrappedBTC
Set-Text "D14" "59.045.06"
$ws.Range("E14").Value = "  -2.17%  "

# Row 15 - Avalanche
Set-Text "D15" "21.05"
$ws.Range("E15").Value = "  -0.99%  "

# Row 16 - ShibaInu
Set-Text "D16" "0.0000136"
$ws.Range("E16").Value = "  -1.52%  "

# Row 17 - WrappedEther
Set-Text "D17" "2.645.19"
$ws.Range("E17").Value = "  -4.12%  "

# Row 18 - BitcoinCash
Set-Text "D18" "338.95"
$ws.Range("E18").Value = "  -3.24%  "

# Row 19 - Polkadot (price unchanged)
$ws.Range("E19").Value = "  -3.08%  "

# Row 20 - Chainlink (price unchanged)
$ws.Range("E20").Value = "  -2.19%  "

# Row 21 - Uniswap
Set-Text "D21" "6.41"
$ws.Range("E21").Value = "  +1.35%  "

# Row 22 - Dai
Set-Text "D22" "0.994"
$ws.Range("E22").Value = "  -0.41%  "

# Row 23 - Litecoin
Set-Text "D23" "64.37"
$ws.Range("E23").Value = "  +2.57%  "

# Row 24 - Polygon
Set-Text "D24" "0.419"
$ws.Range("E24").Value = "  -0.95%  "

# Row 25 - Kaspa (price unchanged)
$ws.Range("E25").Value = "  -1.64%  "

# Row 26 - Binance-PegBSC-USD
Set-Text "D26" "0.998"
$ws.Range("E26").Value = "  +0.33%  "

# Row 27 - PEPE (price unchanged)
$ws.Range("E27").Value = "  -1.44%  "

# Row 28 - InternetComputer(DFINITY)
Set-Text "D28" "7.11"
$ws.Range("E28").Value = "  -1.89%  "

# Row 29 - Aptos (price unchanged)
$ws.Range("E29").Value = "  -2.31%  "

# Row 30 - USDe
Set-Text "D30" "0.998"
$ws.Range("E30").Value = "  +0.01%  "

# Row 31 - PancakeSwap
Set-Text "D31" "1.59"
$ws.Range("E31").Value = "  -0.10%  "

# Row 32 - EthereumClassic (price unchanged)
$ws.Range("E32").Value = "  -1.10%  "

# Row 33 - Monero
Set-Text "D33" "150.75"
$ws.Range("E33").Value = "  +1.87%  "

# Row 34 - NEARProtocol
Set-Text "D34" "4.15"
$ws.Range("E34").Value = "  -3.76%  "

# Row 35 - ImmutableX (price unchanged)
$ws.Range("E35").Value = "  -4.28%  "

# Row 36 - SuiNetwork
Set-Text "D36" "0.899"
$ws.Range("E36").Value = "  -5.63%  "

# Row 37 - Fetch.AI
Set-Text "D37" "0.871"
$ws.Range("E37").Value = "  -0.52%  "

# Row 38 - OKB
Set-Text "D38" "36.95"
$ws.Range("E38").Value = "  +0.31%  "

# Row 39 - Stacks (price unchanged)
$ws.Range("E39").Value = "  -5.49%  "

# Row 40 - Filecoin
Set-Text "D40" "3.58"
$ws.Range("E40").Value = "  -2.82%  "

# Row 41 - Mantle
Set-Text "D41" "0.616"
$ws.Range("E41").Value = "  +0.91%  "

# Row 42 - FirstDigitalUSD (price unchanged)
$ws.Range("E42").Value = "  +0.19%  "

# Row 43 - Bittensor
Set-Text "D43" "276.03"
$ws.Range("E43").Value = "  -1.93%  "

# Row 44 - EnergySwap
Set-Text "D44" "19.84"
$ws.Range("E44").Value = "  -0.46%  "

# Row 45 - Stellar (price unchanged)
$ws.Range("E45").Value = "  -1.68%  "

# Row 46 - WhiteBITCoin
Set-Text "D46" "10.66"
$ws.Range("E46").Value = "  +2.00%  "

# Rows 47-48 swap: Maker/Hedera reorder (Hedera moves up to 47, Maker to 48)
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-Text "D47" "0.0532"
$ws.Range("E47").Value = "  -1.08%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-Text "D48" "2.051.49"
$ws.Range("E48").Value = "  -3.97%  "

# Row 49 - RenderToken (price unchanged)
$ws.Range("E49").Value = "  -3.10%  "

# Row 50 - VeChain (price unchanged)
$ws.Range("E50").Value = "  -2.21%  "

# Row 51 - InjectiveProtocol
Set-Text "D51" "18.89"
$ws.Range("E51").Value = "  -1.52%  "
